$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells (Price/Volume columns, and any Coin/Link text)
# keep their string type instead of being auto-converted to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.99%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.67%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.607"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08348"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.07%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.045"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.01%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9755"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.561"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.89%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1160"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.28%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1920"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.87%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.35"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.43%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09994"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.20%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04663"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.21%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.70%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001287"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006059"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.25%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.004632"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "7.12%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.377"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.41%"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.475"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.38%"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3363"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.20%"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1392"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.94%"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2652"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.17%"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04195"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.95%"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001312"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.84%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.60%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.07%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02765"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "6.98%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05803"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007721"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.12%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.67%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007268"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.96%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.00%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008086"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.15%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3403"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007302"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.72%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.24%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005811"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.01%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003496"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-6.84%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003505"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.24%"
